# Handles float input without breaking stuff
# Updates the marksheet summary numbers/markings and the per-question
# "Student Ans" grid so that previously "float-ish" marked answers are
# now correctly recognised as attempted (and scored), collapses the
# now-unused 3rd answer block (G/H), and drops the D/E block for all
# but the first few questions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Summary block (rows 10-12): give the row-label cells (col A) the same
#    "mtitleStyle" formatting already used by the header row (row 9), then
#    update the numeric figures to their corrected values.
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B10").Value = 22
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 88
$ws.Range("E12").Value = "88/112"

# ---------------------------------------------------------------------------
# 2) The third answer block (columns G/H) is no longer used - clear it out
#    completely (rows 15-21).
# ---------------------------------------------------------------------------
$ws.Range("G15:H21").Clear()

# ---------------------------------------------------------------------------
# 3) The second answer block (columns D/E) is only kept for the first three
#    questions (rows 16-18); remove it for the remaining question rows.
# ---------------------------------------------------------------------------
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------------
# 4) Populate the "Student Ans" cells (column A, and column D for rows
#    16-18) that are now recognised as correctly-answered, giving them the
#    green "correctStyle" formatting (copied from a cell already using it)
#    and the matching option text.
# ---------------------------------------------------------------------------
$ws.Range("B10").Copy()
$correctCells = "A16","D16","D17","A18","D18","A19","A21","A22","A23","A25","A26","A27","A28","A30","A31","A32","A33","A34","A36","A38","A39","A40"
foreach ($addr in $correctCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "Option A"
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("A18").Value = "Option B"
$ws.Range("D18").Value = "Option D"
$ws.Range("A19").Value = "Option C"
$ws.Range("A21").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A23").Value = "Option D"
$ws.Range("A25").Value = "Option A"
$ws.Range("A26").Value = "Option C"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A34").Value = "Option B"
$ws.Range("A36").Value = "Option A"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"
$ws.Range("A40").Value = "Option D"
